# Add new column 'Correction' (column N) to the Card24 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Existing "Event" column (M) was blank for every data row; the author's edit
# filled it in with the literal "nan" placeholder used throughout this sheet.
$ws.Range("M2").Value = "nan"
$ws.Range("M3").Value = "nan"
$ws.Range("M4").Value = "nan"
$ws.Range("M5").Value = "nan"
$ws.Range("M6").Value = "nan"
$ws.Range("M7").Value = "nan"
$ws.Range("M8").Value = "nan"
$ws.Range("M9").Value = "nan"
$ws.Range("M10").Value = "nan"
$ws.Range("M11").Value = "nan"
$ws.Range("M12").Value = "nan"

# New header, matching the same (bold/centered/bordered) style used by the
# other column headers in row 1.
$ws.Range("N1").Value = "Correction"
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

# New "Correction" column cells for each data row start out blank, matching
# the formatting of the other (unstyled) blank data cells on the sheet.
$ws.Range("A2").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N5").PasteSpecial(-4122)
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N10").PasteSpecial(-4122)
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("N12").PasteSpecial(-4122)

$excel.CutCopyMode = 0
